# Week 5 - Deliverable 2 Complete
#
# Logbook.docx edit: adds a new list item to the "Week 4 Workshop" row and
# fills in the next (previously empty) row with "23/04/25" / "Week 5
# Workshop" plus its list of bullet items.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' + $bodyXml + '</w:body>' +
      '</w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- Row 6 ("16/04/25" / "Week 4 Workshop"): add a new bullet to the
#     existing (previously empty) trailing paragraph in cell 2.
$row6Cell2 = $t.Rows.Item(6).Cells.Item(2)
$row6Body = '<w:p><w:pPr>' +
    '<w:pStyle w:val="ListParagraph"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
    '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Further work on deliverable 2</w:t></w:r>' +
  '</w:p>'
$row6Cell2.Range.InsertXML((New-PkgXml $row6Body))

# --- Row 7 cell 1: date "23/04/25"
$row7Cell1 = $t.Rows.Item(7).Cells.Item(1)
$row7Cell1Body = '<w:p><w:pPr>' +
    '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>23/04/25</w:t></w:r>' +
  '</w:p>'
$row7Cell1.Range.InsertXML((New-PkgXml $row7Cell1Body))

# --- Row 7 cell 2: "Week 5 Workshop" + 4 bullet items
$row7Cell2 = $t.Rows.Item(7).Cells.Item(2)
$row7Cell2Body =
  '<w:p><w:pPr>' +
      '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Week 5 Workshop</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Further work on deliverable 2</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:rPr><w:bCs/><w:color w:val="FF0000"/><w:lang w:val="en-AU"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Investigated </w:t></w:r>' +
    '<w:r><w:rPr><w:bCs/><w:color w:val="FF0000"/><w:lang w:val="en-AU"/></w:rPr><w:t>AlAinSat-1</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Worked together with Thomas Unipan. Had him go over my design notes</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr>' +
      '<w:pStyle w:val="ListParagraph"/>' +
      '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>Investigated Greek cubesat</w:t></w:r>' +
  '</w:p>'
$row7Cell2.Range.InsertXML((New-PkgXml $row7Cell2Body))

Write-Output "Done applying Week 5 logbook edits."
